$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.084461850806488
$ws.Range("D2").Value = 1.035966942793223
$ws.Range("E2").Value = 1.084941446035002
$ws.Range("F2").Value = 1.091366756559735
$ws.Range("I2").Value = 1.035196290025188
$ws.Range("J2").Value = 1.089321515017736
$ws.Range("K2").Value = 1.038762032748761
$ws.Range("L2").Value = 1.087603161880338
$ws.Range("M2").Value = 1.094011932972122
$ws.Range("N2").Value = 1.090868476515985
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.088400404273401
$ws.Range("D3").Value = 1.036543744345057
$ws.Range("E3").Value = 1.088429186497859
$ws.Range("F3").Value = 1.094844316924935
$ws.Range("I3").Value = 1.035421127384873
$ws.Range("J3").Value = 1.092909335410998
$ws.Range("K3").Value = 1.039148956674079
$ws.Range("L3").Value = 1.090903600364728
$ws.Range("M3").Value = 1.097303484020031
$ws.Range("N3").Value = 1.094461392025734
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.090924434630436
$ws.Range("D4").Value = 1.036914510511331
$ws.Range("E4").Value = 1.090663065324709
$ws.Range("F4").Value = 1.09707130875045
$ws.Range("I4").Value = 1.035562981794735
$ws.Range("J4").Value = 1.095207038626339
$ws.Range("K4").Value = 1.039396164703685
$ws.Range("L4").Value = 1.093016223342038
$ws.Range("M4").Value = 1.099410046338128
$ws.Range("N4").Value = 1.096762358243192
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.091979871746557
$ws.Range("D5").Value = 1.037069801690718
$ws.Range("E5").Value = 1.091596872910764
$ws.Range("F5").Value = 1.09800214749844
$ws.Range("I5").Value = 1.035621759432255
$ws.Range("J5").Value = 1.096167461697796
$ws.Range("K5").Value = 1.039499345993129
$ws.Range("L5").Value = 1.09389903528981
$ws.Range("M5").Value = 1.100290234520502
$ws.Range("N5").Value = 1.097724145225572
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.09215675762349
$ws.Range("D6").Value = 1.037095842120045
$ws.Range("E6").Value = 1.09175335656075
$ws.Range("F6").Value = 1.098158128304331
$ws.Range("I6").Value = 1.035631578463281
$ws.Range("J6").Value = 1.096328401730321
$ws.Range("K6").Value = 1.039516627212677
$ws.Range("L6").Value = 1.094046955234576
$ws.Range("M6").Value = 1.100437709426954
$ws.Range("N6").Value = 1.097885313811404
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.090938559438754
$ws.Range("D7").Value = 1.036916587780225
$ws.Range("E7").Value = 1.090675563566089
$ws.Range("F7").Value = 1.097083767613324
$ws.Range("I7").Value = 1.035563770540896
$ws.Range("J7").Value = 1.095219893345451
$ws.Range("K7").Value = 1.039397546329859
$ws.Range("L7").Value = 1.093028040261818
$ws.Range("M7").Value = 1.099421828499279
$ws.Range("N7").Value = 1.096775231217479
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.085798090569221
$ws.Range("D8").Value = 1.036162392148509
$ws.Range("E8").Value = 1.086124996075454
$ws.Range("F8").Value = 1.092546927181295
$ws.Range("I8").Value = 1.03527303340394
$ws.Range("J8").Value = 1.09053908541338
$ws.Range("K8").Value = 1.038893456038625
$ws.Range("L8").Value = 1.088723421144775
$ws.Range("M8").Value = 1.095129251827845
$ws.Range("N8").Value = 1.092087776001248
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.076543443515051
$ws.Range("D9").Value = 1.034814037402121
$ws.Range("E9").Value = 1.077922757168827
$ws.Range("F9").Value = 1.084366662922773
$ws.Range("I9").Value = 1.034732393033917
$ws.Range("J9").Value = 1.082099852988376
$ws.Range("K9").Value = 1.037980488802018
$ws.Range("L9").Value = 1.080954380381415
$ws.Range("M9").Value = 1.087379135381537
$ws.Range("N9").Value = 1.083636558898208
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.070228543173357
$ws.Range("D10").Value = 1.033901380400824
$ws.Range("E10").Value = 1.072319627274391
$ws.Range("F10").Value = 1.078776812129345
$ws.Range("I10").Value = 1.03435216714612
$ws.Range("J10").Value = 1.076333197543704
$ws.Range("K10").Value = 1.03735449389085
$ws.Range("L10").Value = 1.07564028882697
$ws.Range("M10").Value = 1.082076227721406
$ws.Range("N10").Value = 1.07786171414136
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.067456666291284
$ws.Range("D11").Value = 1.033502754621064
$ws.Range("E11").Value = 1.06985870888646
$ws.Range("F11").Value = 1.076321346016968
$ws.Range("I11").Value = 1.034182652213086
$ws.Range("J11").Value = 1.073800044451019
$ws.Range("K11").Value = 1.037079140636121
$ws.Range("L11").Value = 1.073304669805365
$ws.Range("M11").Value = 1.079745137520465
$ws.Range("N11").Value = 1.075324963680726
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.066421166417073
$ws.Range("D12").Value = 1.033354154947185
$ws.Range("E12").Value = 1.068939157629157
$ws.Range("F12").Value = 1.075403777683133
$ws.Range("I12").Value = 1.034118938598193
$ws.Range("J12").Value = 1.072853437369337
$ws.Range("K12").Value = 1.036976201384216
$ws.Range("L12").Value = 1.072431689528706
$ws.Range("M12").Value = 1.078873794574731
$ws.Range("N12").Value = 1.074377012308411
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.066643556167392
$ws.Range("D13").Value = 1.03338605442865
$ws.Range("E13").Value = 1.069136655374655
$ws.Range("F13").Value = 1.075600851958568
$ws.Range("I13").Value = 1.034132639506124
$ws.Range("J13").Value = 1.073056749015882
$ws.Range("K13").Value = 1.036998312366822
$ws.Range("L13").Value = 1.072619196193082
$ws.Range("M13").Value = 1.079060952018984
$ws.Range("N13").Value = 1.074580612680814
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.06737119353593
$ws.Range("D14").Value = 1.033490482277436
$ws.Range("E14").Value = 1.069782811225134
$ws.Range("F14").Value = 1.076245613113097
$ws.Range("I14").Value = 1.034177400983498
$ws.Range("J14").Value = 1.073721914990178
$ws.Range("K14").Value = 1.03707064521915
$ws.Range("L14").Value = 1.07323262107476
$ws.Range("M14").Value = 1.079673225034027
$ws.Range("N14").Value = 1.075246723267089
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.067818724677491
$ws.Range("D15").Value = 1.033554752706263
$ws.Range("E15").Value = 1.070180198788813
$ws.Range("F15").Value = 1.076642135828546
$ws.Range("I15").Value = 1.034204880374052
$ws.Range("J15").Value = 1.074130985237815
$ws.Range("K15").Value = 1.037115123831389
$ws.Range("L15").Value = 1.07360984611246
$ws.Range("M15").Value = 1.08004973449229
$ws.Range("N15").Value = 1.075656374441399
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.070411691000217
$ws.Range("D16").Value = 1.033927762104891
$ws.Range("E16").Value = 1.072482198478592
$ws.Range("F16").Value = 1.078939015417183
$ws.Range("I16").Value = 1.034363313253227
$ws.Range("J16").Value = 1.076500531797173
$ws.Range("K16").Value = 1.037372676394541
$ws.Range("L16").Value = 1.075794547997888
$ws.Range("M16").Value = 1.082230179908862
$ws.Range("N16").Value = 1.078029286028662
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.07202796789994
$ws.Range("D17").Value = 1.034160809474046
$ws.Range("E17").Value = 1.073916719030672
$ws.Range("F17").Value = 1.080370246251287
$ws.Range("I17").Value = 1.034461377959338
$ws.Range("J17").Value = 1.077977033564314
$ws.Range("K17").Value = 1.037533071961177
$ws.Range("L17").Value = 1.077155534523563
$ws.Range("M17").Value = 1.083588414009665
$ws.Range("N17").Value = 1.079507884597669
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.072967117078462
$ws.Range("D18").Value = 1.034296411313858
$ws.Range("E18").Value = 1.074750116705691
$ws.Range("F18").Value = 1.081201696200107
$ws.Range("I18").Value = 1.034518108221606
$ws.Range("J18").Value = 1.078834781195583
$ws.Range("K18").Value = 1.037626214840484
$ws.Range("L18").Value = 1.077946054396245
$ws.Range("M18").Value = 1.08437729799562
$ws.Range("N18").Value = 1.08036685032896
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.073286740439619
$ws.Range("D19").Value = 1.034342592397649
$ws.Range("E19").Value = 1.075033725392906
$ws.Range("F19").Value = 1.081484635790963
$ws.Range("I19").Value = 1.034537372656925
$ws.Range("J19").Value = 1.079126669658527
$ws.Range("K19").Value = 1.037657904546174
$ws.Range("L19").Value = 1.078215044899672
$ws.Range("M19").Value = 1.08464572547245
$ws.Range("N19").Value = 1.080659153307001
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.07185493073642
$ws.Range("D20").Value = 1.034135840018171
$ws.Range("E20").Value = 1.073763155090212
$ws.Range("F20").Value = 1.080217038237419
$ws.Range("I20").Value = 1.034450905190743
$ws.Range("J20").Value = 1.077818979677434
$ws.Range("K20").Value = 1.037515905879724
$ws.Range("L20").Value = 1.077009858653249
$ws.Range("M20").Value = 1.083443036581283
$ws.Range("N20").Value = 1.079349606256141
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.067157087599632
$ws.Range("D21").Value = 1.033459745702191
$ws.Range("E21").Value = 1.069592687034669
$ws.Range("F21").Value = 1.076055900773328
$ws.Range("I21").Value = 1.034164240629631
$ws.Range("J21").Value = 1.073526199077733
$ws.Range("K21").Value = 1.037049363368445
$ws.Range("L21").Value = 1.073052134469644
$ws.Range("M21").Value = 1.079493078844853
$ws.Range("N21").Value = 1.0750507294156
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.064169081486994
$ws.Range("D22").Value = 1.033031570314304
$ws.Range("E22").Value = 1.066938852073649
$ws.Range("F22").Value = 1.073407688904795
$ws.Range("I22").Value = 1.033979665888499
$ws.Range("J22").Value = 1.07079415626793
$ws.Range("K22").Value = 1.036752198085251
$ws.Range("L22").Value = 1.070532233058336
$ws.Range("M22").Value = 1.076977802080197
$ws.Range("N22").Value = 1.072314806791637
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.06575642309962
$ws.Range("D23").Value = 1.033258852305796
$ws.Range("E23").Value = 1.068348787229182
$ws.Range("F23").Value = 1.074814665387155
$ws.Range("I23").Value = 1.034077929066181
$ws.Range("J23").Value = 1.072245678131476
$ws.Range("K23").Value = 1.036910099662714
$ws.Range("L23").Value = 1.071871148482992
$ws.Range("M23").Value = 1.078314289599432
$ws.Range("N23").Value = 1.073768389982721
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.071933129882132
$ws.Range("D24").Value = 1.034147123661731
$ws.Range("E24").Value = 1.073832554311196
$ws.Range("F24").Value = 1.080286276719955
$ws.Range("I24").Value = 1.034455638833027
$ws.Range("J24").Value = 1.077890408134508
$ws.Range("K24").Value = 1.037523663768772
$ws.Range("L24").Value = 1.07707569355409
$ws.Range("M24").Value = 1.083508736723312
$ws.Range("N24").Value = 1.079421136149817
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.078960586154284
$ws.Range("D25").Value = 1.035164987776235
$ws.Range("E25").Value = 1.080066147066145
$ws.Range("F25").Value = 1.086504622373916
$ws.Range("I25").Value = 1.034875593896271
$ws.Range("J25").Value = 1.084305448770311
$ws.Range("K25").Value = 1.038219514060957
$ws.Range("L25").Value = 1.082985763919559
$ws.Range("M25").Value = 1.089405887414619
$ws.Range("N25").Value = 1.085845286879137
